# Applies the pedalboard-soundcard BoM update:
#  - R10 is no longer fitted: drop it from the "R5 R6 R7 R10" group on the
#    BoM sheet (now "R5 R6 R7", qty/PCB 4 -> 3) and add it as its own row on
#    the DNF (Do Not Fit) sheet.
#  - Component/summary counters drop by one fitted part (19 -> 20 groups is
#    unaffected-looking only because the group count actually grows by one
#    new DNF-only group; Fitted/Total Components go 72 -> 71).

$wb = $excel.ActiveWorkbook
$wsBom = $wb.Worksheets.Item("BoM")
$wsDnf = $wb.Worksheets.Item("DNF")

# ---------------------------------------------------------------------
# Helper: write a numeric-looking string value into a cell while keeping
# it a genuine text cell (so it still matches t="s" shared-string cells
# like the "Row"/"Quantity Per PCB" columns use), instead of letting the
# engine coerce "3" -> the number 3. A plain `.Value = "3"` assignment
# gets auto-coerced to a number, and forcing NumberFormat="@" first stamps
# a brand new (quote-prefixed) style onto the cell - so instead we stage
# the text in a scratch cell formatted as Text, copy it, and paste-special
# *values only* onto the target, which keeps the target's original style.
# ---------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $scratch = $range.Worksheet.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

# --- BoM sheet -----------------------------------------------------------

# Component Groups: 19 -> 20
$wsBom.Range("F2").Value = 20

# Fitted Components: "72 (69 SMD/ 3 THT)" -> "71 (68 SMD/ 3 THT)"
$wsBom.Range("F4").Value = "71 (68 SMD/ 3 THT)"

# Total Components: 72 -> 71
$wsBom.Range("F6").Value = 71

# R10 drops out of the R5/R6/R7/R10 3K9 group
$wsBom.Range("D18").Value = "R5 R6 R7"

# Quantity Per PCB for that group: 4 -> 3 (kept as text, like the source)
Set-TextValue $wsBom.Range("G18") "3"

# --- DNF sheet -------------------------------------------------------------

# Component Groups: 19 -> 20
$wsDnf.Range("F2").Value = 20

# Fitted Components: "72 (69 SMD/ 3 THT)" -> "71 (68 SMD/ 3 THT)"
$wsDnf.Range("F4").Value = "71 (68 SMD/ 3 THT)"

# Total Components: 72 -> 71
$wsDnf.Range("F6").Value = 71

# New row 10: R10, now DNF. Clone formatting (incl. the blank B column)
# from BoM!A10:J10, which uses the same alternating-stripe style set.
$wsBom.Range("A10:J10").Copy($wsDnf.Range("A10:J10"))
$wsDnf.Rows.Item(10).RowHeight = 30

$wsDnf.Range("C10").Value = "R"
$wsDnf.Range("D10").Value = "R10"
$wsDnf.Range("E10").Value = "3K9"
$wsDnf.Range("F10").Value = "R_0805_2012Metric"
$wsDnf.Range("H10").Value = " (DNF)"
$wsDnf.Range("I10").Value = "https://www.seielect.com/catalog/sei-rmcf_rmcp.pdf"
$wsDnf.Range("J10").Value = "https://www.digikey.ch/de/products/detail/stackpole-electronics-inc/RMCF0805FT3K90/1760599"

# A10 ("Row" = 2) and G10 ("Quantity Per PCB" = 1) stay text cells
Set-TextValue $wsDnf.Range("A10") "2"
Set-TextValue $wsDnf.Range("G10") "1"
